$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '23.446.08'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +1.02%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.641.06'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +2.33%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '305.17'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3732'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.33%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '52.42'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.17%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.3620'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -0.25%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.254'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("E11").Value = '  -0.02%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '22.83'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.02%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.600'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.06%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.00001269'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.82%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '7.282'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.89%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.630.87'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.74%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '94.34'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.41%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06863'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -0.13%  '

$ws.Range("E21").Value = '  -0.46%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.9996'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '23.446.50'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '12.73'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.76%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '3.037'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.83%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.406'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.79%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '21.22'
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '151.74'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.16%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '5.297'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.95%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '135.52'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.22%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.292'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -3.35%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.811.82'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.81%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '6.756'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.9535'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -1.33%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.02838'
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '10.34'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +1.02%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.2519'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.30%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.07228'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -3.84%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.08782'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.20%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '6.054'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '1.376'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +0.40%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.7050'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '12.46'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.68%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '16.10'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +3.22%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.6514'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("E46").Value = '  +0.38%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.9997'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '4.009'
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.07968'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '128.31'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.06%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.198'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.64%  '
